$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1773049645390071
$ws.Range("C2").Value = 0.5638297872340425
$ws.Range("J2").Value = 0.02836879432624113
$ws.Range("P2").Value = 0.1418439716312057
$ws.Range("S2").Value = 0.08865248226950355
$ws.Range("B3").Value = 0.005847953216374269
$ws.Range("C3").Value = 0.01169590643274854
$ws.Range("J3").Value = 0.03508771929824561
$ws.Range("P3").Value = 0.7777777777777778
$ws.Range("S3").Value = 0.1695906432748538
$ws.Range("B6").Value = 0.0625
$ws.Range("F6").Value = 0.05833333333333333
$ws.Range("J6").Value = 0.2458333333333333
$ws.Range("O6").Value = 0.05416666666666667
$ws.Range("Q6").Value = 0.1166666666666667
$ws.Range("R6").Value = 0.09583333333333334
$ws.Range("S6").Value = 0.3666666666666666
$ws.Range("B7").Value = 0.09523809523809523
$ws.Range("D7").Value = 0.009523809523809525
$ws.Range("E7").Value = 0.004761904761904762
$ws.Range("F7").Value = 0.06666666666666667
$ws.Range("J7").Value = 0.1571428571428571
$ws.Range("O7").Value = 0.02380952380952381
$ws.Range("Q7").Value = 0.2047619047619048
$ws.Range("R7").Value = 0.08571428571428572
$ws.Range("S7").Value = 0.3523809523809524
$ws.Range("B8").Value = 0.09195402298850575
$ws.Range("D8").Value = 0.01839080459770115
$ws.Range("F8").Value = 0.05057471264367816
$ws.Range("J8").Value = 0.1402298850574713
$ws.Range("O8").Value = 0.02068965517241379
$ws.Range("Q8").Value = 0.1471264367816092
$ws.Range("R8").Value = 0.135632183908046
$ws.Range("S8").Value = 0.3954022988505747
$ws.Range("B9").Value = 0.095
$ws.Range("D9").Value = 0.01
$ws.Range("F9").Value = 0.035
$ws.Range("J9").Value = 0.145
$ws.Range("O9").Value = 0.02
$ws.Range("Q9").Value = 0.175
$ws.Range("R9").Value = 0.115
$ws.Range("S9").Value = 0.405
$ws.Range("B10").Value = 0.09314079422382672
$ws.Range("D10").Value = 0.01516245487364621
$ws.Range("E10").Value = 0.0007220216606498195
$ws.Range("F10").Value = 0.0779783393501805
$ws.Range("J10").Value = 0.1220216606498195
$ws.Range("O10").Value = 0.02382671480144404
$ws.Range("Q10").Value = 0.2245487364620939
$ws.Range("R10").Value = 0.1003610108303249
$ws.Range("S10").Value = 0.3422382671480144
$ws.Range("G11").Value = 0.1420454545454546
$ws.Range("J11").Value = 0.1107954545454545
$ws.Range("K11").Value = 0.1903409090909091
$ws.Range("L11").Value = 0.5511363636363636
$ws.Range("S11").Value = 0.005681818181818182
$ws.Range("G12").Value = 0.719047619047619
$ws.Range("J12").Value = 0.1761904761904762
$ws.Range("K12").Value = 0.01428571428571429
$ws.Range("L12").Value = 0.06190476190476191
$ws.Range("S12").Value = 0.02857142857142857
$ws.Range("G13").Value = 0.5277777777777778
$ws.Range("J13").Value = 0.3888888888888889
$ws.Range("S13").Value = 0.08333333333333333
$ws.Range("F15").Value = 0.01879699248120301
$ws.Range("H15").Value = 0.1090225563909774
$ws.Range("I15").Value = 0.06015037593984962
$ws.Range("J15").Value = 0.3796992481203008
$ws.Range("K15").Value = 0.06766917293233082
$ws.Range("M15").Value = 0.007518796992481203
$ws.Range("O15").Value = 0.06390977443609022
$ws.Range("S15").Value = 0.2932330827067669
$ws.Range("F16").Value = 0.03370786516853932
$ws.Range("H16").Value = 0.1573033707865168
$ws.Range("I16").Value = 0.101123595505618
$ws.Range("J16").Value = 0.4269662921348314
$ws.Range("K16").Value = 0.101123595505618
$ws.Range("M16").Value = 0.02808988764044944
$ws.Range("O16").Value = 0.03932584269662921
$ws.Range("S16").Value = 0.1123595505617977
$ws.Range("F17").Value = 0.01659751037344398
$ws.Range("H17").Value = 0.1929460580912863
$ws.Range("I17").Value = 0.0975103734439834
$ws.Range("J17").Value = 0.4128630705394191
$ws.Range("K17").Value = 0.0954356846473029
$ws.Range("M17").Value = 0.01037344398340249
$ws.Range("N17").Value = 0.002074688796680498
$ws.Range("O17").Value = 0.06846473029045644
$ws.Range("S17").Value = 0.1037344398340249
$ws.Range("F18").Value = 0.02298850574712644
$ws.Range("H18").Value = 0.1877394636015326
$ws.Range("I18").Value = 0.07662835249042145
$ws.Range("J18").Value = 0.4137931034482759
$ws.Range("K18").Value = 0.1226053639846743
$ws.Range("M18").Value = 0.01532567049808429
$ws.Range("O18").Value = 0.08045977011494253
$ws.Range("S18").Value = 0.08045977011494253
$ws.Range("F19").Value = 0.02019386106623586
$ws.Range("H19").Value = 0.1873990306946688
$ws.Range("I19").Value = 0.08239095315024232
$ws.Range("J19").Value = 0.3796445880452343
$ws.Range("K19").Value = 0.130048465266559
$ws.Range("M19").Value = 0.01857835218093699
$ws.Range("N19").Value = 0.001615508885298869
$ws.Range("O19").Value = 0.07835218093699516
$ws.Range("S19").Value = 0.1017770597738288
